$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Trait lists" header in G1 is replaced with "Tags" (the trait-import
# template dropped the old "Trait lists" column and added a new "Tags"
# column in its place).
$ws.Range("G1").Value = "Tags"

# Move/leave the selection on the edited cell, matching the saved view state.
$ws.Range("G1").Select()

